$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2921
$ws.Range("B2").Value = 6910
$ws.Range("C2").Value = 3480
$ws.Range("D2").Value = 6600
$ws.Range("E2").Value = 4155
$ws.Range("F2").Value = 7921
